# Vehicle-Load data update:
# "reading and posting entities to the database implemented"
#
# - "Vehicles" sheet: Fuel Type "Diesel" -> "Gasoline" (rows 2 and 3)
# - "Loads" sheet: Load Type "Locker" -> "Trailer" (rows 4 and 5)
# - "Loads" sheet: normalise the font on F2:G3 so it matches the rest of
#   the column (F4:G5 already used this font) - this is what collapses the
#   stray duplicate font/style when the sheet is re-saved.
# - Leave the cursor/selection where the author left it when the file was
#   saved.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Vehicles")
$ws2 = $wb.Worksheets.Item("Loads")

# --- Vehicles: Fuel Type column (I) ---
$ws1.Range("I2").Value = "Gasoline"
$ws1.Range("I3").Value = "Gasoline"

# --- Loads: Load Type column (J) ---
$ws2.Range("J4").Value = "Trailer"
$ws2.Range("J5").Value = "Trailer"

# --- Loads: normalise F2:G3 font to match the rest of the column ---
$fmtRange = $ws2.Range("F2:G3")
$fmtRange.Font.Name = "Arial"
$fmtRange.Font.Size = 11

# --- restore cursor/selection position on each sheet ---
$ws1.Range("F20").Select()
$ws2.Range("D5").Select()
